$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.497.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +9.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.614.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +9.40%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.87%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9901'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3689'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3420'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +11.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.28'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +7.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07087'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9994'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.82'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +9.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.938'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.674'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.30%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.607.33'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +9.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9905'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06760'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +13.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.35'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +12.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.048'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +10.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.13'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +11.51%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.507.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +9.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.386'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.564'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +20.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.07'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.61'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +13.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.786.85'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.86'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.08%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.209'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +23.70%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.043'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9565'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +18.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08259'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.667'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.77%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.279'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.57%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '12.01'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +15.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.275'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.624'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +16.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06130'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02233'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +8.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2028'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9906'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5944'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +11.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.830'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +8.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.14'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5729'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.39'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.984'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +9.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06824'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.09'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.78%  '
